$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 211, shifting existing rows 211+ down by one.
$ws.Rows.Item(211).Insert()

# Populate the newly inserted row 211 with the new price observation.
$ws.Cells.Item(211, 1).Value = 4
$ws.Cells.Item(211, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(211, 3).Value = "Los Lagos"
$ws.Cells.Item(211, 4).Value = 44754
$ws.Cells.Item(211, 4).NumberFormat = $ws.Cells.Item(212, 4).NumberFormat
$ws.Cells.Item(211, 5).Value = 10
$ws.Cells.Item(211, 6).Value = "Fruta"
$ws.Cells.Item(211, 7).Value = 100108
$ws.Cells.Item(211, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(211, 9).Value = 100108005
$ws.Cells.Item(211, 10).Value = "Piña"
$ws.Cells.Item(211, 11).Value = "Caramelo"
$ws.Cells.Item(211, 12).Value = "Primera"
$ws.Cells.Item(211, 13).Value = 160
$ws.Cells.Item(211, 14).Value = 22000
$ws.Cells.Item(211, 15).Value = 23000
$ws.Cells.Item(211, 16).Value = 22500
$ws.Cells.Item(211, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(211, 18).Value = "Ecuador"
$ws.Cells.Item(211, 19).Value = 1875
$ws.Cells.Item(211, 20).Value = 12
